$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.697917
$ws.Range("H2").Value = 2.093751
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.837667
$ws.Range("N2").Value = 32.513001
$ws.Range("O2").Value = 0.5477430134675739
$ws.Range("P2").Value = 0.5477430134675739
$ws.Range("Q2").Value = 7.563792039639001
$ws.Range("R2").Value = 68.074128356751
$ws.Range("S2").Value = 0.5477430134675739
$ws.Range("T2").Value = 0.5477430134675739

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.697917
$ws.Range("H3").Value = 2.093751
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.948376333333334
$ws.Range("N3").Value = 26.845129
$ws.Range("O3").Value = 0.4522569865324261
$ws.Range("P3").Value = 0.4522569865324261
$ws.Range("Q3").Value = 6.245223965431
$ws.Range("R3").Value = 56.207015688879
$ws.Range("S3").Value = 0.4522569865324261
$ws.Range("T3").Value = 0.4522569865324261
